# Auto-generated script to update Leve price/profit columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3950.6667
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H17").Value = 1959.2667
$ws.Range("J17").Value = 1959.2667
$ws.Range("L17").Value = 5877.800099999999
$ws.Range("N17").Value = -6213.800099999999
$ws.Range("H40").Value = 2239.8
$ws.Range("I40").Value = 1979.6
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1979.6
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -1804.6
$ws.Range("N40").Value = -2850
$ws.Range("H43").Value = 9749.5
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
$ws.Range("H62").Value = 5808.1665
$ws.Range("I62").Value = 3514
$ws.Range("K62").Value = 3514
$ws.Range("M62").Value = -2890
$ws.Range("H65").Value = 5808.1665
$ws.Range("I65").Value = 3514
$ws.Range("K65").Value = 17570
$ws.Range("M65").Value = -14450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 390
$ws.Range("I4").Value = 383.33334
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 383.33334
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -267.33334
$ws.Range("N4").Value = -632
$ws.Range("H61").Value = 6662
$ws.Range("I61").Value = 6662
$ws.Range("K61").Value = 6662
$ws.Range("M61").Value = -6450
$ws.Range("H74").Value = 2488.3845
$ws.Range("I74").Value = 2084.9
$ws.Range("K74").Value = 2084.9
$ws.Range("M74").Value = -1210.9
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 2488.3845
$ws.Range("I77").Value = 2084.9
$ws.Range("K77").Value = 10424.5
$ws.Range("M77").Value = -6056.5
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H122").Value = 3928.7778
$ws.Range("I122").Value = 3071.8
$ws.Range("K122").Value = 9215.400000000001
$ws.Range("M122").Value = -6765.400000000001
$ws.Range("H136").Value = 6662
$ws.Range("I136").Value = 6662
$ws.Range("K136").Value = 19986
$ws.Range("M136").Value = -17436

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 773392.4
$ws.Range("J20").Value = 1668132.6
$ws.Range("L20").Value = 1668132.6
$ws.Range("N20").Value = -1668626.6
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H86").Value = 1430144
$ws.Range("I86").Value = 1252.25
$ws.Range("K86").Value = 1252.25
$ws.Range("M86").Value = -129.25
$ws.Range("H89").Value = 1430144
$ws.Range("I89").Value = 1252.25
$ws.Range("K89").Value = 6261.25
$ws.Range("M89").Value = -645.25
$ws.Range("H94").Value = 97.666664
$ws.Range("I94").Value = 97.666664
$ws.Range("K94").Value = 97.666664
$ws.Range("M94").Value = 353.333336
$ws.Range("H107").Value = 1249.591
$ws.Range("I107").Value = 932.35
$ws.Range("K107").Value = 932.35
$ws.Range("M107").Value = 987.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1043.4117
$ws.Range("I16").Value = 1133.3636
$ws.Range("J16").Value = 878.5
$ws.Range("K16").Value = 1133.3636
$ws.Range("L16").Value = 878.5
$ws.Range("M16").Value = -846.3635999999999
$ws.Range("N16").Value = -1452.5
$ws.Range("H22").Value = 148246.28
$ws.Range("I22").Value = 340333
$ws.Range("K22").Value = 340333
$ws.Range("M22").Value = -339983
$ws.Range("H31").Value = 4922.933
$ws.Range("I31").Value = 3290.1
$ws.Range("J31").Value = 8188.6
$ws.Range("K31").Value = 3290.1
$ws.Range("L31").Value = 8188.6
$ws.Range("M31").Value = -2995.1
$ws.Range("N31").Value = -8778.6
$ws.Range("H34").Value = 4922.933
$ws.Range("I34").Value = 3290.1
$ws.Range("J34").Value = 8188.6
$ws.Range("K34").Value = 3290.1
$ws.Range("L34").Value = 8188.6
$ws.Range("M34").Value = -3088.1
$ws.Range("N34").Value = -8592.6
$ws.Range("H41").Value = 14665
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14665
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14665
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -15521
$ws.Range("H113").Value = 1043.4117
$ws.Range("I113").Value = 1133.3636
$ws.Range("J113").Value = 878.5
$ws.Range("K113").Value = 1133.3636
$ws.Range("L113").Value = 878.5
$ws.Range("M113").Value = 1036.6364
$ws.Range("N113").Value = -5218.5
$ws.Range("H132").Value = 4215.923
$ws.Range("I132").Value = 2613.5
$ws.Range("J132").Value = 6779.8
$ws.Range("K132").Value = 7840.5
$ws.Range("L132").Value = 20339.4
$ws.Range("M132").Value = -5310.5
$ws.Range("N132").Value = -25399.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1689.2858
$ws.Range("I11").Value = 1275
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 3825
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = -3685
$ws.Range("N11").Value = -6280
$ws.Range("H92").Value = 724.25
$ws.Range("J92").Value = 765.6667
$ws.Range("L92").Value = 2297.0001
$ws.Range("N92").Value = -4793.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2240
$ws.Range("I5").Value = 2175
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 2175
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = -2063
$ws.Range("N5").Value = -2724
$ws.Range("H74").Value = 55000
$ws.Range("J74").Value = 55000
$ws.Range("L74").Value = 55000
$ws.Range("N74").Value = -56872
$ws.Range("H77").Value = 55000
$ws.Range("J77").Value = 55000
$ws.Range("L77").Value = 165000
$ws.Range("N77").Value = -174360
$ws.Range("H80").Value = 2826.3044
$ws.Range("I80").Value = 2826.3044
$ws.Range("K80").Value = 2826.3044
$ws.Range("M80").Value = -1828.3044
$ws.Range("H83").Value = 2826.3044
$ws.Range("I83").Value = 2826.3044
$ws.Range("K83").Value = 14131.522
$ws.Range("M83").Value = -9139.522000000001
$ws.Range("H102").Value = 2631.7646
$ws.Range("I102").Value = 934
$ws.Range("K102").Value = 934
$ws.Range("M102").Value = 688
$ws.Range("H122").Value = 42009.54
$ws.Range("I122").Value = 2870.0667
$ws.Range("J122").Value = 95381.55
$ws.Range("K122").Value = 8610.2001
$ws.Range("L122").Value = 286144.65
$ws.Range("M122").Value = -6160.2001
$ws.Range("N122").Value = -291044.65
$ws.Range("H132").Value = 4503.1113
$ws.Range("I132").Value = 2750
$ws.Range("K132").Value = 8250
$ws.Range("M132").Value = -5720

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 6000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 6000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -6224
$ws.Range("H22").Value = 1111.125
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 1127.1428
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 1127.1428
$ws.Range("M22").Value = -704
$ws.Range("N22").Value = -1717.1428
$ws.Range("H24").Value = 50020004
$ws.Range("I24").Value = 100000000
$ws.Range("J24").Value = 40007
$ws.Range("K24").Value = 100000000
$ws.Range("L24").Value = 40007
$ws.Range("M24").Value = -99999657
$ws.Range("N24").Value = -40693
$ws.Range("H27").Value = 1111.125
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 1127.1428
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 1127.1428
$ws.Range("M27").Value = -892
$ws.Range("N27").Value = -1341.1428
$ws.Range("H93").Value = 7316
$ws.Range("I93").Value = 7316
$ws.Range("K93").Value = 7316
$ws.Range("M93").Value = -6068
$ws.Range("H103").Value = 10392.8
$ws.Range("J103").Value = 10392.8
$ws.Range("L103").Value = 10392.8
$ws.Range("N103").Value = -12736.8
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 18000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3100400.2
$ws.Range("J4").Value = 1250500.2
$ws.Range("L4").Value = 1250500.2
$ws.Range("N4").Value = -1250726.2
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H54").Value = 24083.25
$ws.Range("J54").Value = 19999.834
$ws.Range("L54").Value = 19999.834
$ws.Range("N54").Value = -21039.834
$ws.Range("H107").Value = 3018.818
$ws.Range("I107").Value = 2401
$ws.Range("J107").Value = 4666.3335
$ws.Range("K107").Value = 7203
$ws.Range("L107").Value = 13999.0005
$ws.Range("M107").Value = -5283
$ws.Range("N107").Value = -17839.0005
$ws.Range("H113").Value = 353.15384
$ws.Range("I113").Value = 281
$ws.Range("K113").Value = 843
$ws.Range("M113").Value = 1327
$ws.Range("H122").Value = 1323.7368
$ws.Range("I122").Value = 1355.8889
$ws.Range("J122").Value = 745
$ws.Range("K122").Value = 4067.6667
$ws.Range("L122").Value = 2235
$ws.Range("M122").Value = -1617.6667
$ws.Range("N122").Value = -7135
$ws.Range("H135").Value = 59000
$ws.Range("J135").Value = 59000
$ws.Range("L135").Value = 59000
$ws.Range("N135").Value = -69140

